$d = $word.ActiveDocument

# --- 1. Insert the two new sub-bullet paragraphs after
#        "Absolute value of left and right looks" ------------------------
$lastPara = $d.Paragraphs.Last
$insPoint = $lastPara.Range
$insPoint.Collapse(0)
$null = $insPoint.InsertParagraphAfter()

$para2 = $d.Paragraphs.Last
$para2.Range.ListFormat.ListLevelNumber = 2
$body2 = $para2.Range.Duplicate
$null = $body2.MoveEnd(1, -1)
$body2.Text = "Moving Average of absolute value of distances from median"

$endPoint = $d.Paragraphs.Last.Range
$endPoint.Collapse(0)
$null = $endPoint.InsertParagraphAfter()

$para3 = $d.Paragraphs.Last
$para3.Range.ListFormat.ListLevelNumber = 2
$body3 = $para3.Range.Duplicate
$null = $body3.MoveEnd(1, -1)
# A temporary trailing marker character is used so the later bookmark
# insertion point falls *inside* the run rather than exactly on the
# run/paragraph-mark boundary (inserting exactly at that boundary makes
# the host snap the zero-length bookmark around the whole run instead of
# leaving it collapsed after the text).
$body3.Text = "Moving Average of occurrences (1 or 0) X"

# --- 2. Re-anchor the "_GoBack" bookmark to the end of the new last
#        paragraph (this also removes it from its old location, since a
#        document can only have one bookmark with a given name). --------
$body3Again = $d.Paragraphs.Last.Range.Duplicate
$null = $body3Again.MoveEnd(1, -1)
$bmPos = $body3Again.End - 1
$bmRange = $d.Paragraphs.Last.Range.Duplicate
$null = $bmRange.SetRange($bmPos, $bmPos)
$d.Bookmarks.Add("_GoBack", $bmRange)

# --- 3. Remove the temporary trailing "X" marker used above. -----------
$trail = $d.Paragraphs.Last.Range.Duplicate
$null = $trail.MoveEnd(1, -1)
$null = $trail.MoveStart(1, $trail.End - $trail.Start - 1)
$trail.Text = ""
